$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25/26: Polygon and PEPE swap positions, with updated price/volume values
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.504"
$ws.Range("E25").Value = "  -6.22%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -5.91%  "

# Remaining per-row price/volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.233.97"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.174.83"
$ws.Range("E3").Value = "  -8.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.22"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.81"
$ws.Range("E6").Value = "  -5.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.175.36"
$ws.Range("E9").Value = "  -8.47%  "
$ws.Range("E10").Value = "  -7.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -5.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.726.49"
$ws.Range("E13").Value = "  -8.50%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.33"
$ws.Range("E15").Value = "  -9.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.235.94"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("E17").Value = "  -5.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.175.55"
$ws.Range("E18").Value = "  -8.32%  "
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").Value = "  -6.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.89"
$ws.Range("E21").Value = "  -5.35%  "
$ws.Range("E22").Value = "  -6.75%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.61"
$ws.Range("E24").Value = "  -6.53%  "
$ws.Range("E27").Value = "  -5.13%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.52"
$ws.Range("E31").Value = "  -7.88%  "
$ws.Range("E32").Value = "  -5.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.98"
$ws.Range("E33").Value = "  -7.47%  "
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("E36").Value = "  -7.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.78"
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.818"
$ws.Range("E38").Value = "  -7.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.92"
$ws.Range("E39").Value = "  -8.12%  "
$ws.Range("E40").Value = "  -7.06%  "
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.614.70"
$ws.Range("E42").Value = "  -7.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.17"
$ws.Range("E43").Value = "  -7.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.38"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("E45").Value = "  -8.36%  "
$ws.Range("E46").Value = "  -6.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.66"
$ws.Range("E47").Value = "  -6.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "318.88"
$ws.Range("E48").Value = "  -6.43%  "
$ws.Range("E49").Value = "  -7.84%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.01%  "
